$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded, so insert a new row at
# position 18 (shifting the existing rows 18-78 down to 19-79) and
# populate it with the new data point.
$ws.Rows.Item(18).Insert(-4121)

$ws.Cells.Item(18, 1).Value = 9
$ws.Cells.Item(18, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(18, 3).Value = "Metropolitana"
$ws.Cells.Item(18, 4).Value = 45133
$ws.Cells.Item(18, 5).Value = 13
$ws.Cells.Item(18, 6).Value = 100112010
$ws.Cells.Item(18, 7).Value = "Achicoria"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 70
$ws.Cells.Item(18, 11).Value = 7000
$ws.Cells.Item(18, 12).Value = 7000
$ws.Cells.Item(18, 13).Value = 7000
$ws.Cells.Item(18, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(18, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(18, 16).Value = 438
$ws.Cells.Item(18, 17).Value = 16
$ws.Cells.Item(18, 18).Value = "Hortaliza"
